# Auto-generated edit script: updates hard-coded market-data cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ================= Sheet: ALC =================
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 89.5
$ws.Range("I11").Value = 89.5
$ws.Range("K11").Value = 89.5
$ws.Range("M11").Value = 50.5

# Row 29
$ws.Range("H29").Value = 1740.6666
$ws.Range("I29").Value = 450
$ws.Range("J29").Value = 1998.8
$ws.Range("K29").Value = 1350
$ws.Range("L29").Value = 5996.4
$ws.Range("M29").Value = -1069
$ws.Range("N29").Value = -6558.4

# Row 38
$ws.Range("H38").Value = 93.333336
$ws.Range("I38").Value = 93.333336
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 280.000008
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 91.99999200000002
$ws.Range("N38").ClearContents()

# Row 87
$ws.Range("H87").Value = 29483.482
$ws.Range("J87").Value = 29483.482
$ws.Range("L87").Value = 29483.482
$ws.Range("N87").Value = -31979.482

# Row 90
$ws.Range("H90").Value = 29483.482
$ws.Range("J90").Value = 29483.482
$ws.Range("L90").Value = 88450.446
$ws.Range("N90").Value = -100930.446

# Row 139
$ws.Range("H139").Value = 51875.555
$ws.Range("J139").Value = 51875.555
$ws.Range("L139").Value = 51875.555
$ws.Range("N139").Value = -62155.555

# ================= Sheet: ARM =================
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3507.72
$ws.Range("I32").Value = 2694.5054
$ws.Range("J32").Value = 18958.8
$ws.Range("K32").Value = 2694.5054
$ws.Range("L32").Value = 18958.8
$ws.Range("M32").Value = -2407.5054
$ws.Range("N32").Value = -19532.8

# Row 42
$ws.Range("H42").Value = 15000
$ws.Range("J42").Value = 15000
$ws.Range("L42").Value = 15000
$ws.Range("N42").Value = -15972

# Row 110
$ws.Range("H110").Value = 21444.2
$ws.Range("I110").Value = 30566
$ws.Range("J110").Value = 2060.375
$ws.Range("K110").Value = 30566
$ws.Range("L110").Value = 2060.375
$ws.Range("M110").Value = -28521
$ws.Range("N110").Value = -6150.375

# Row 132
$ws.Range("H132").Value = 2771.4443
$ws.Range("I132").Value = 1657.2778
$ws.Range("J132").Value = 4999.778
$ws.Range("K132").Value = 4971.8334
$ws.Range("L132").Value = 14999.334
$ws.Range("M132").Value = -2441.8334
$ws.Range("N132").Value = -20059.334

# Row 133
$ws.Range("H133").Value = 36837.832
$ws.Range("J133").Value = 42205.4
$ws.Range("L133").Value = 42205.4
$ws.Range("N133").Value = -47265.4

# ================= Sheet: BSM =================
$ws = $wb.Worksheets.Item("BSM")
# Row 45
$ws.Range("H45").Value = 45950
$ws.Range("J45").Value = 45950
$ws.Range("L45").Value = 45950
$ws.Range("N45").Value = -47566

# ================= Sheet: CRP =================
$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 5800
$ws.Range("I25").Value = 5800
$ws.Range("K25").Value = 5800
$ws.Range("M25").Value = -5626

# Row 51
$ws.Range("H51").Value = 31755.75
$ws.Range("J51").Value = 31755.75
$ws.Range("L51").Value = 31755.75
$ws.Range("N51").Value = -33227.75

# Row 61
$ws.Range("H61").Value = 31755.75
$ws.Range("J61").Value = 31755.75
$ws.Range("L61").Value = 31755.75
$ws.Range("N61").Value = -32451.75

# Row 132
$ws.Range("H132").Value = 1604.415
$ws.Range("I132").Value = 578.2632
$ws.Range("J132").Value = 4204
$ws.Range("K132").Value = 1734.7896
$ws.Range("L132").Value = 12612
$ws.Range("M132").Value = 795.2103999999999
$ws.Range("N132").Value = -17672

# Row 134
$ws.Range("H134").Value = 1731.0625
$ws.Range("I134").Value = 1779.8
$ws.Range("K134").Value = 5339.4
$ws.Range("M134").Value = -2804.4

# ================= Sheet: CUL =================
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 9980
$ws.Range("I3").Value = 9980
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 29940
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -29828
$ws.Range("N3").ClearContents()

# Row 112
$ws.Range("H112").Value = 1361.1724
$ws.Range("I112").Value = 1168
$ws.Range("J112").Value = 1383.4615
$ws.Range("K112").Value = 3504
$ws.Range("L112").Value = 4150.3845
$ws.Range("M112").Value = -2396
$ws.Range("N112").Value = -6366.3845

# Row 131
$ws.Range("H131").Value = 5330
$ws.Range("J131").Value = 5965.909
$ws.Range("L131").Value = 17897.727
$ws.Range("N131").Value = -27977.727

# Row 137
$ws.Range("H137").Value = 2026.8334
$ws.Range("I137").Value = 2404.875
$ws.Range("J137").Value = 1889.3636
$ws.Range("K137").Value = 7214.625
$ws.Range("L137").Value = 5668.0908
$ws.Range("M137").Value = -2114.625
$ws.Range("N137").Value = -15868.0908

# ================= Sheet: GSM =================
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2584.9565
$ws.Range("I102").Value = 2265.45
$ws.Range("K102").Value = 2265.45
$ws.Range("M102").Value = -643.4499999999998

# Row 122
$ws.Range("H122").Value = 2487.8928
$ws.Range("I122").Value = 2205.25
$ws.Range("J122").Value = 3194.5
$ws.Range("K122").Value = 6615.75
$ws.Range("L122").Value = 9583.5
$ws.Range("M122").Value = -4165.75
$ws.Range("N122").Value = -14483.5

# Row 126
$ws.Range("H126").Value = 1924.6471
$ws.Range("I126").Value = 1721.9
$ws.Range("J126").Value = 2214.2856
$ws.Range("K126").Value = 5165.700000000001
$ws.Range("L126").Value = 6642.8568
$ws.Range("M126").Value = -2695.700000000001
$ws.Range("N126").Value = -11582.8568

# Row 132
$ws.Range("H132").Value = 2043.6923
$ws.Range("I132").Value = 1771.2174
$ws.Range("J132").Value = 4132.6665
$ws.Range("K132").Value = 5313.6522
$ws.Range("L132").Value = 12397.9995
$ws.Range("M132").Value = -2783.6522
$ws.Range("N132").Value = -17457.9995

# ================= Sheet: LTW =================
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 7243.9272
$ws.Range("I132").Value = 8570.767
$ws.Range("J132").Value = 5651.72
$ws.Range("K132").Value = 25712.301
$ws.Range("L132").Value = 16955.16
$ws.Range("M132").Value = -23182.301
$ws.Range("N132").Value = -22015.16

# Row 140
$ws.Range("H140").Value = 69925.5
$ws.Range("J140").Value = 69925.5
$ws.Range("L140").Value = 69925.5
$ws.Range("N140").Value = -80285.5

# ================= Sheet: WVR =================
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1388.5
$ws.Range("I96").Value = 880
$ws.Range("J96").Value = 1490.2
$ws.Range("K96").Value = 880
$ws.Range("L96").Value = 1490.2
$ws.Range("M96").Value = 493
$ws.Range("N96").Value = -4236.2

# Row 132
$ws.Range("H132").Value = 3291.0625
$ws.Range("I132").Value = 3047.0715
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 9141.2145
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6611.2145
$ws.Range("N132").Value = -20057

